$d = $word.ActiveDocument

# 1) "редактировать любые посты, кроме опубликованных (свои и писателей);"
#    -> "редактировать любые посты, кроме опубликованных и скрытых (свои и писателей);"
$d.Content.Find.Execute(
    "редактировать любые посты, кроме опубликованных (свои и писателей);",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "редактировать любые посты, кроме опубликованных и скрытых (свои и писателей);", 2)

# 2) "Модератор может редактировать любые посты, кроме уже опубликованных, публиковать ..."
#    -> "... кроме уже опубликованных и скрытых, публиковать ..."
$d.Content.Find.Execute(
    "Модератор может редактировать любые посты, кроме уже опубликованных, публиковать",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Модератор может редактировать любые посты, кроме уже опубликованных и скрытых, публиковать", 2)
